# Update Excel file with latest predictions
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Home win" - rows shift up (oldest match dropped) + new match added
# ---------------------------------------------------------------------
$wsHome = $wb.Worksheets.Item("Home win")

$wsHome.Range("A2").Value2 = "04-02-2025 20:45"
$wsHome.Range("B2").Value2 = "ENGLAND"
$wsHome.Range("C2").Value2 = "LEAGUE TWO"
$wsHome.Range("D2").Value2 = "Salford City - Bromley"
$wsHome.Range("E2").Value2 = 73.3
$wsHome.Range("F2").Value2 = 2

$wsHome.Range("A3").Value2 = "04-02-2025 20:45"
$wsHome.Range("B3").Value2 = "SCOTLAND"
$wsHome.Range("C3").Value2 = "LEAGUE TWO"
$wsHome.Range("D3").Value2 = "Elgin City - Bonnyrigg Rose Athletic"
$wsHome.Range("E3").Value2 = 73.3
$wsHome.Range("F3").Value2 = 2

$wsHome.Range("A4").Value2 = "04-02-2025 20:00"
$wsHome.Range("B4").Value2 = "ENGLAND"
$wsHome.Range("C4").Value2 = "PREMIER LEAGUE CUP"
$wsHome.Range("D4").Value2 = "Ipswich Town U21 - Watford U21"
$wsHome.Range("E4").Value2 = 73.3
$wsHome.Range("F4").Value2 = 1.73

$wsHome.Range("A5").Value2 = "04-02-2025 11:00"
$wsHome.Range("B5").Value2 = "TURKEY"
$wsHome.Range("C5").Value2 = "CUP"
$wsHome.Range("D5").Value2 = "Kocaelispor - Sivasspor"
$wsHome.Range("E5").Value2 = 70
$wsHome.Range("F5").Value2 = 2.25

# ---------------------------------------------------------------------
# Sheet "Btts" - 4 new matches appended (rows 6-9)
# ---------------------------------------------------------------------
$wsBtts = $wb.Worksheets.Item("Btts")

$wsBtts.Range("A6").Value2 = "05-02-2025 21:00"
$wsBtts.Range("B6").Value2 = "ITALY"
$wsBtts.Range("C6").Value2 = "COPPA ITALIA"
$wsBtts.Range("D6").Value2 = "AC Milan - AS Roma"
$wsBtts.Range("E6").Value2 = 80
$wsBtts.Range("F6").Value2 = 1.7

$wsBtts.Range("A7").Value2 = "05-02-2025 21:00"
$wsBtts.Range("B7").Value2 = "NETHERLANDS"
$wsBtts.Range("C7").Value2 = "KNVB BEKER"
$wsBtts.Range("D7").Value2 = "GO Ahead Eagles - Noordwijk"
$wsBtts.Range("E7").Value2 = 84
$wsBtts.Range("F7").Value2 = 2.05

$wsBtts.Range("A8").Value2 = "05-02-2025 20:45"
$wsBtts.Range("B8").Value2 = "FRANCE"
$wsBtts.Range("C8").Value2 = "COUPE DE FRANCE"
$wsBtts.Range("D8").Value2 = "Toulouse - Guingamp"
$wsBtts.Range("E8").Value2 = 83.3
$wsBtts.Range("F8").Value2 = 1.77

$wsBtts.Range("A9").Value2 = "05-02-2025 20:45"
$wsBtts.Range("B9").Value2 = "FRANCE"
$wsBtts.Range("C9").Value2 = "COUPE DE FRANCE"
$wsBtts.Range("D9").Value2 = "Cannes - Dives-Cabourg"
$wsBtts.Range("E9").Value2 = 84
$wsBtts.Range("F9").Value2 = 1.91

# ---------------------------------------------------------------------
# Sheet "Over_Under" - new Slovakia match inserted before the Switzerland
# row (old rows 12-13 shift down to 13-14), plus 2 new matches appended
# (rows 15-16)
# ---------------------------------------------------------------------
$wsOU = $wb.Worksheets.Item("Over_Under")

# Old row 13 (World / Slavia Praha II) moves to row 14 - write it first so
# we never lose data while shuffling rows manually.
$wsOU.Range("A14").Value2 = "04-02-2025 10:30"
$wsOU.Range("B14").Value2 = "WORLD"
$wsOU.Range("C14").Value2 = "FRIENDLIES CLUBS"
$wsOU.Range("D14").Value2 = "Slavia Praha II - Příbram"
$wsOU.Range("E14").Value2 = 80
$wsOU.Range("F14").Value2 = 1.57
$wsOU.Range("G14").Value2 = 53.3
$wsOU.Range("H14").Value2 = 2.4

# Old row 12 (Switzerland / FC Sion - Servette FC) moves to row 13.
$wsOU.Range("A13").Value2 = "04-02-2025 20:30"
$wsOU.Range("B13").Value2 = "SWITZERLAND"
$wsOU.Range("C13").Value2 = "SUPER LEAGUE"
$wsOU.Range("D13").Value2 = "FC Sion - Servette FC"
$wsOU.Range("E13").Value2 = 70
$wsOU.Range("F13").Value2 = 1.83
$wsOU.Range("G13").Value2 = 55
$wsOU.Range("H13").Value2 = 3.1

# New row 12 (Slovakia / Slovan Bratislava - AS Trencin).
$wsOU.Range("A12").Value2 = "04-02-2025 18:00"
$wsOU.Range("B12").Value2 = "SLOVAKIA"
$wsOU.Range("C12").Value2 = "CUP"
$wsOU.Range("D12").Value2 = "Slovan Bratislava - AS Trencin"
$wsOU.Range("E12").Value2 = 75
$wsOU.Range("F12").Value2 = 1.5
$wsOU.Range("G12").Value2 = 55
$wsOU.Range("H12").Value2 = 2.2

# New row 15 (France / Cannes - Dives-Cabourg).
$wsOU.Range("A15").Value2 = "05-02-2025 20:45"
$wsOU.Range("B15").Value2 = "FRANCE"
$wsOU.Range("C15").Value2 = "COUPE DE FRANCE"
$wsOU.Range("D15").Value2 = "Cannes - Dives-Cabourg"
$wsOU.Range("E15").Value2 = 86.7
$wsOU.Range("F15").Value2 = 1.83
$wsOU.Range("G15").Value2 = 40
$wsOU.Range("H15").Value2 = 3

# New row 16 (Switzerland / FC ST. Gallen - FC Lugano).
$wsOU.Range("A16").Value2 = "05-02-2025 20:30"
$wsOU.Range("B16").Value2 = "SWITZERLAND"
$wsOU.Range("C16").Value2 = "SUPER LEAGUE"
$wsOU.Range("D16").Value2 = "FC ST. Gallen - FC Lugano"
$wsOU.Range("E16").Value2 = 75
$wsOU.Range("F16").Value2 = 1.65
$wsOU.Range("G16").Value2 = 55
$wsOU.Range("H16").Value2 = 2.5

Write-Host "Predictions updated"
